$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.178715
$ws.Range("H2").Value = 0.536145
$ws.Range("I2").Value = 0.09904930989061336
$ws.Range("J2").Value = 0.09904930989061336
$ws.Range("M2").Value = 30.10959799999999
$ws.Range("N2").Value = 90.32879399999999
$ws.Range("O2").Value = 0.5062005690267993
$ws.Range("P2").Value = 0.5062005690267992
$ws.Range("Q2").Value = 5.381036806569998
$ws.Range("R2").Value = 48.42933125913
$ws.Range("S2").Value = 0.05013881702834027
$ws.Range("T2").Value = 0.05013881702834026
$ws.Range("G3").Value = 0.178715
$ws.Range("H3").Value = 0.536145
$ws.Range("I3").Value = 0.09904930989061336
$ws.Range("J3").Value = 0.09904930989061336
$ws.Range("O3").Value = 0.2331091635681292
$ws.Range("P3").Value = 0.2331091635681292
$ws.Range("Q3").Value = 2.47800786064
$ws.Range("R3").Value = 22.30207074576
$ws.Range("S3").Value = 0.02308930178060131
$ws.Range("T3").Value = 0.02308930178060131
$ws.Range("G4").Value = 0.178715
$ws.Range("H4").Value = 0.536145
$ws.Range("I4").Value = 0.09904930989061336
$ws.Range("J4").Value = 0.09904930989061336
$ws.Range("M4").Value = 7.300670666666666
$ws.Range("N4").Value = 21.902012
$ws.Range("O4").Value = 0.122738392114831
$ws.Range("P4").Value = 0.1227383921148309
$ws.Range("Q4").Value = 1.304739358193333
$ws.Range("R4").Value = 11.74265422374
$ws.Range("S4").Value = 0.01215715303605751
$ws.Range("T4").Value = 0.0121571530360575
$ws.Range("G5").Value = 0.178715
$ws.Range("H5").Value = 0.536145
$ws.Range("I5").Value = 0.09904930989061336
$ws.Range("J5").Value = 0.09904930989061336
$ws.Range("M5").Value = 0.9177576666666667
$ws.Range("N5").Value = 2.753273
$ws.Range("O5").Value = 0.01542928115796745
$ws.Range("P5").Value = 0.01542928115796744
$ws.Range("Q5").Value = 0.1640170613983333
$ws.Range("R5").Value = 1.476153552585
$ws.Range("S5").Value = 0.00152825965080492
$ws.Range("T5").Value = 0.001528259650804919
$ws.Range("G6").Value = 0.178715
$ws.Range("H6").Value = 0.536145
$ws.Range("I6").Value = 0.09904930989061336
$ws.Range("J6").Value = 0.09904930989061336
$ws.Range("M6").Value = 1.836782666666667
$ws.Range("N6").Value = 5.510348
$ws.Range("O6").Value = 0.03087986864006716
$ws.Range("P6").Value = 0.03087986864006715
$ws.Range("Q6").Value = 0.3282606142733334
$ws.Range("R6").Value = 2.95434552846
$ws.Range("S6").Value = 0.003058629678311445
$ws.Range("T6").Value = 0.003058629678311444
$ws.Range("G7").Value = 0.178715
$ws.Range("H7").Value = 0.536145
$ws.Range("I7").Value = 0.09904930989061336
$ws.Range("J7").Value = 0.09904930989061336
$ws.Range("M7").Value = 5.451052
$ws.Range("N7").Value = 16.353156
$ws.Range("O7").Value = 0.09164272549220594
$ws.Range("P7").Value = 0.09164272549220591
$ws.Range("Q7").Value = 0.9741847581799998
$ws.Range("R7").Value = 8.767662823619998
$ws.Range("S7").Value = 0.009077148716497918
$ws.Range("T7").Value = 0.009077148716497917
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 1.450498333333333
$ws.Range("H8").Value = 4.351495
$ws.Range("I8").Value = 0.8039104659046613
$ws.Range("J8").Value = 0.8039104659046612
$ws.Range("M8").Value = 30.10959799999999
$ws.Range("N8").Value = 90.32879399999999
$ws.Range("O8").Value = 0.5062005690267993
$ws.Range("P8").Value = 0.5062005690267992
$ws.Range("Q8").Value = 43.67392171633666
$ws.Range("R8").Value = 393.0652954470299
$ws.Range("S8").Value = 0.4069399352875389
$ws.Range("T8").Value = 0.4069399352875387
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 1.450498333333333
$ws.Range("H9").Value = 4.351495
$ws.Range("I9").Value = 0.8039104659046613
$ws.Range("J9").Value = 0.8039104659046612
$ws.Range("O9").Value = 0.2331091635681292
$ws.Range("P9").Value = 0.2331091635681292
$ws.Range("Q9").Value = 20.11216893850666
$ws.Range("R9").Value = 181.00952044656
$ws.Range("S9").Value = 0.1873988962907007
$ws.Range("T9").Value = 0.1873988962907006
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 1.450498333333333
$ws.Range("H10").Value = 4.351495
$ws.Range("I10").Value = 0.8039104659046613
$ws.Range("J10").Value = 0.8039104659046612
$ws.Range("M10").Value = 7.300670666666666
$ws.Range("N10").Value = 21.902012
$ws.Range("O10").Value = 0.122738392114831
$ws.Range("P10").Value = 0.1227383921148309
$ws.Range("Q10").Value = 10.58961063421555
$ws.Range("R10").Value = 95.30649570793999
$ws.Range("S10").Value = 0.09867067798942276
$ws.Range("T10").Value = 0.09867067798942271
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 1.450498333333333
$ws.Range("H11").Value = 4.351495
$ws.Range("I11").Value = 0.8039104659046613
$ws.Range("J11").Value = 0.8039104659046612
$ws.Range("M11").Value = 0.9177576666666667
$ws.Range("N11").Value = 2.753273
$ws.Range("O11").Value = 0.01542928115796745
$ws.Range("P11").Value = 0.01542928115796744
$ws.Range("Q11").Value = 1.331205965903889
$ws.Range("R11").Value = 11.980853693135
$ws.Range("S11").Value = 0.01240376060427562
$ws.Range("T11").Value = 0.01240376060427562
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 1.450498333333333
$ws.Range("H12").Value = 4.351495
$ws.Range("I12").Value = 0.8039104659046613
$ws.Range("J12").Value = 0.8039104659046612
$ws.Range("M12").Value = 1.836782666666667
$ws.Range("N12").Value = 5.510348
$ws.Range("O12").Value = 0.03087986864006716
$ws.Range("P12").Value = 0.03087986864006715
$ws.Range("Q12").Value = 2.664250196695556
$ws.Range("R12").Value = 23.97825177026
$ws.Range("S12").Value = 0.02482464958551113
$ws.Range("T12").Value = 0.02482464958551112
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 1.450498333333333
$ws.Range("H13").Value = 4.351495
$ws.Range("I13").Value = 0.8039104659046613
$ws.Range("J13").Value = 0.8039104659046612
$ws.Range("M13").Value = 5.451052
$ws.Range("N13").Value = 16.353156
$ws.Range("O13").Value = 0.09164272549220594
$ws.Range("P13").Value = 0.09164272549220591
$ws.Range("Q13").Value = 7.906741840913332
$ws.Range("R13").Value = 71.16067656822
$ws.Range("S13").Value = 0.07367254614721226
$ws.Range("T13").Value = 0.07367254614721222
$ws.Range("E14").Value = 1
$ws.Range("F14").Value = 0.3333333333333333
$ws.Range("G14").Value = 0.17509
$ws.Range("H14").Value = 0.52527
$ws.Range("I14").Value = 0.09704022420472538
$ws.Range("J14").Value = 0.09704022420472537
$ws.Range("M14").Value = 30.10959799999999
$ws.Range("N14").Value = 90.32879399999999
$ws.Range("O14").Value = 0.5062005690267993
$ws.Range("P14").Value = 0.5062005690267992
$ws.Range("Q14").Value = 5.271889513819999
$ws.Range("R14").Value = 47.44700562438
$ws.Range("S14").Value = 0.04912181671092018
$ws.Range("T14").Value = 0.04912181671092015
$ws.Range("E15").Value = 1
$ws.Range("F15").Value = 0.3333333333333333
$ws.Range("G15").Value = 0.17509
$ws.Range("H15").Value = 0.52527
$ws.Range("I15").Value = 0.09704022420472538
$ws.Range("J15").Value = 0.09704022420472537
$ws.Range("O15").Value = 0.2331091635681292
$ws.Range("P15").Value = 0.2331091635681292
$ws.Range("Q15").Value = 2.42774471264
$ws.Range("R15").Value = 21.84970241376
$ws.Range("S15").Value = 0.02262096549682726
$ws.Range("T15").Value = 0.02262096549682725
$ws.Range("E16").Value = 1
$ws.Range("F16").Value = 0.3333333333333333
$ws.Range("G16").Value = 0.17509
$ws.Range("H16").Value = 0.52527
$ws.Range("I16").Value = 0.09704022420472538
$ws.Range("J16").Value = 0.09704022420472537
$ws.Range("M16").Value = 7.300670666666666
$ws.Range("N16").Value = 21.902012
$ws.Range("O16").Value = 0.122738392114831
$ws.Range("P16").Value = 0.1227383921148309
$ws.Range("Q16").Value = 1.278274427026667
$ws.Range("R16").Value = 11.50446984324
$ws.Range("S16").Value = 0.01191056108935069
$ws.Range("T16").Value = 0.01191056108935069
$ws.Range("E17").Value = 1
$ws.Range("F17").Value = 0.3333333333333333
$ws.Range("G17").Value = 0.17509
$ws.Range("H17").Value = 0.52527
$ws.Range("I17").Value = 0.09704022420472538
$ws.Range("J17").Value = 0.09704022420472537
$ws.Range("M17").Value = 0.9177576666666667
$ws.Range("N17").Value = 2.753273
$ws.Range("O17").Value = 0.01542928115796745
$ws.Range("P17").Value = 0.01542928115796744
$ws.Range("Q17").Value = 0.1606901898566667
$ws.Range("R17").Value = 1.44621170871
$ws.Range("S17").Value = 0.001497260902886906
$ws.Range("T17").Value = 0.001497260902886905
$ws.Range("E18").Value = 1
$ws.Range("F18").Value = 0.3333333333333333
$ws.Range("G18").Value = 0.17509
$ws.Range("H18").Value = 0.52527
$ws.Range("I18").Value = 0.09704022420472538
$ws.Range("J18").Value = 0.09704022420472537
$ws.Range("M18").Value = 1.836782666666667
$ws.Range("N18").Value = 5.510348
$ws.Range("O18").Value = 0.03087986864006716
$ws.Range("P18").Value = 0.03087986864006715
$ws.Range("Q18").Value = 0.3216022771066667
$ws.Range("R18").Value = 2.89442049396
$ws.Range("S18").Value = 0.002996589376244585
$ws.Range("T18").Value = 0.002996589376244584
$ws.Range("E19").Value = 1
$ws.Range("F19").Value = 0.3333333333333333
$ws.Range("G19").Value = 0.17509
$ws.Range("H19").Value = 0.52527
$ws.Range("I19").Value = 0.09704022420472538
$ws.Range("J19").Value = 0.09704022420472537
$ws.Range("M19").Value = 5.451052
$ws.Range("N19").Value = 16.353156
$ws.Range("O19").Value = 0.09164272549220594
$ws.Range("P19").Value = 0.09164272549220591
$ws.Range("Q19").Value = 0.9544246946799999
$ws.Range("R19").Value = 8.589822252119999
$ws.Range("S19").Value = 0.008893030628495767
$ws.Range("T19").Value = 0.008893030628495762
